$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.588.10"
$ws.Range("E2").Value = "  +7.56%  "
$ws.Range("D3").Value = "3.635.75"
$ws.Range("E3").Value = "  +7.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("D8").Value = "3.607.19"
$ws.Range("E8").Value = "  +6.95%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.76%  "
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000298"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("D15").Value = "4.212.80"
$ws.Range("E15").Value = "  +7.60%  "
$ws.Range("D16").Value = "3.630.90"
$ws.Range("E16").Value = "  +7.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.58%  "
$ws.Range("D18").Value = "70.457.07"
$ws.Range("E18").Value = "  +7.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.85%  "
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("E21").Value = "  +5.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "496.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +17.28%  "
$ws.Range("E25").Value = "  +8.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.92%  "
$ws.Range("E28").Value = "  +6.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "617.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.24%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.117"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.29%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.51%  "
$ws.Range("D36").Value = "0.0₃0831"
$ws.Range("E36").Value = "  +12.86%  "
$ws.Range("E37").Value = "  +5.22%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.53%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.403"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.51%  "
$ws.Range("D42").Value = "3.352.36"
$ws.Range("E42").Value = "  +7.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0446"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.31%  "
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.65%  "
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
